# Insert two new rows (451, 452) into the Betarraga price sheet for the
# new weekly observation (date serial 44918 = 2022-12-23), pushing all
# existing rows 451..499 down to 453..501.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A451:A452").EntireRow.Insert()

# New row 451: "Primera" quality observation
$ws.Cells.Item(451, 1).Value = 7
$ws.Cells.Item(451, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(451, 3).Value = "Ñuble"
$ws.Cells.Item(451, 4).Value = 44918
$ws.Cells.Item(451, 5).Value = 16
$ws.Cells.Item(451, 6).Value = 100114014
$ws.Cells.Item(451, 7).Value = "Betarraga"
$ws.Cells.Item(451, 8).Value = "Sin especificar"
$ws.Cells.Item(451, 9).Value = "Primera"
$ws.Cells.Item(451, 10).Value = 400
$ws.Cells.Item(451, 11).Value = 700
$ws.Cells.Item(451, 12).Value = 800
$ws.Cells.Item(451, 13).Value = 750
$ws.Cells.Item(451, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(451, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(451, 16).Value = 150
$ws.Cells.Item(451, 17).Value = 5
$ws.Cells.Item(451, 18).Value = "Hortaliza"

# New row 452: "Segunda" quality observation, same date
$ws.Cells.Item(452, 1).Value = 7
$ws.Cells.Item(452, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(452, 3).Value = "Ñuble"
$ws.Cells.Item(452, 4).Value = 44918
$ws.Cells.Item(452, 5).Value = 16
$ws.Cells.Item(452, 6).Value = 100114014
$ws.Cells.Item(452, 7).Value = "Betarraga"
$ws.Cells.Item(452, 8).Value = "Sin especificar"
$ws.Cells.Item(452, 9).Value = "Segunda"
$ws.Cells.Item(452, 10).Value = 300
$ws.Cells.Item(452, 11).Value = 600
$ws.Cells.Item(452, 12).Value = 600
$ws.Cells.Item(452, 13).Value = 600
$ws.Cells.Item(452, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(452, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(452, 16).Value = 120
$ws.Cells.Item(452, 17).Value = 5
$ws.Cells.Item(452, 18).Value = "Hortaliza"
